$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "71.894.92"
$ws.Range("E2").Value = "  -0.61%  "

# Row 3
$ws.Range("D3").Value = "3.993.83"
$ws.Range("E3").Value = "  -1.00%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.18%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "543.29"
$ws.Range("E5").Value = "  +4.52%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.64"
$ws.Range("E6").Value = "  +1.57%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.701"
$ws.Range("E7").Value = "  +12.52%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.01%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.744"
$ws.Range("E9").Value = "  +1.31%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.169"
$ws.Range("E10").Value = "  -3.15%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "51.49"

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000323"
$ws.Range("E12").Value = "  -3.21%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.64"
$ws.Range("E13").Value = "  -2.39%  "

# Row 14
$ws.Range("D14").Value = "4.636.24"
$ws.Range("E14").Value = "  -1.16%  "

# Row 15
$ws.Range("D15").Value = "3.996.30"
$ws.Range("E15").Value = "  -0.90%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.04"
$ws.Range("E16").Value = "  -0.68%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.43"
$ws.Range("E17").Value = "  -3.62%  "

# Row 18
$ws.Range("E18").Value = "  -0.27%  "

# Row 19
$ws.Range("E19").Value = "  -2.12%  "

# Row 20
$ws.Range("D20").Value = "71.760.80"
$ws.Range("E20").Value = "  -0.79%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "428.73"
$ws.Range("E21").Value = "  -1.79%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "97.10"
$ws.Range("E22").Value = "  -1.42%  "

# Row 23
$ws.Range("E23").Value = "  -0.99%  "

# Row 24
$ws.Range("E24").Value = "  +5.84%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.28"
$ws.Range("E25").Value = "  -2.71%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.34"
$ws.Range("E26").Value = "  -4.58%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.69"
$ws.Range("E27").Value = "  -5.07%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.85"
$ws.Range("E28").Value = "  +1.09%  "

# Row 29
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "36.68"
$ws.Range("E29").Value = "  -1.62%  "

# Row 30
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.59"
$ws.Range("E30").Value = "  +16.80%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.39"
$ws.Range("E31").Value = "  -0.77%  "

# Row 32
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.130"
$ws.Range("E32").Value = "  +1.42%  "

# Row 33
$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "49.13"
$ws.Range("E33").Value = "  +20.33%  "

# Row 34
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.21"
$ws.Range("E34").Value = "  +3.29%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "676.11"
$ws.Range("E35").Value = "  -2.06%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "65.77"
$ws.Range("E36").Value = "  -3.47%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.442"
$ws.Range("E37").Value = "  +0.83%  "

# Row 38
$ws.Range("D38").Value = "0.0₃0831"
$ws.Range("E38").Value = "  -6.85%  "

# Row 39
$ws.Range("E39").Value = "  -2.23%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.40"
$ws.Range("E40").Value = "  -7.82%  "

# Row 41
$ws.Range("E41").Value = "  +0.01%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.33"
$ws.Range("E42").Value = "  +5.57%  "

# Row 43
$ws.Range("E43").Value = "  +0.15%  "

# Row 44
$ws.Range("E44").Value = "  -0.64%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.149"
$ws.Range("E45").Value = "  +2.63%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.70"
$ws.Range("E46").Value = "  -2.98%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.80"
$ws.Range("E47").Value = "  +8.49%  "

# Row 48
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.33"
$ws.Range("E48").Value = "  -5.02%  "

# Row 49
$ws.Range("B49").Value = "FLOKI"
$ws.Range("C49").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000282"
$ws.Range("E49").Value = "  +3.18%  "

# Row 50
$ws.Range("E50").Value = "  -3.79%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "145.00"
$ws.Range("E51").Value = "  +1.77%  "
